$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 8 corresponds to CNZMonumentSketchMaxSize; column F holds the Korean
# translation string. Update its text to the new translation.
$ws.Range("F8").Value = "기념비의 비석 최대 개수 (최소 3)"
